$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "45.449.27"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -6.22%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.681.09"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.26%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.34%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.38"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.24%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.57"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -7.26%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.597"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -3.00%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.34%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.581"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.83%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.55"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -5.42%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0848"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.12%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.10"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.91%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.085.46"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.39%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.70%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.679.45"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.27%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.934"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.41%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "15.17"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.83%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "45.580.09"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -5.54%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.02%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.90"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.19%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.82"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.60%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "74.87"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.41%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "282.70"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.49%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.06"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.88%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.24"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.49%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "30.95"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.77%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.14%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.05"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.65%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.63"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.49%  "

$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.22"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -4.37%  "

$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "38.39"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -5.13%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.21"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.84%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.78"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.25%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.35"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +4.22%  "

$ws.Range("B35").Value = "Monero"
$ws.Range("C35").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "155.21"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.26%  "

$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0846"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.28%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.82"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.77%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.30%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "25.78"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +12.82%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.124"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.89%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "16.31"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.96%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.64"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.33%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0327"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.44%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.00"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -8.03%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.109.78"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -4.22%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.18%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "93.27"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -5.61%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "111.87"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.91%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.29"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -6.82%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.937.06"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.54%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.200"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.86%  "
